$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 gets a new date, hours and description entry
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = 42429

$ws.Range("B12").Value = 3

$ws.Range("C12").Value = "Thema aangepast, bugfixes en nieuws carousel slider"

# Move the active selection from C17 to C18
$ws.Range("C18").Select()
